$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had three always-blank helper columns (E, F and I) in the
# imported-file block (B:M were used, but only B-D/G-P actually carried
# data). Remove those blank columns outright so the real data shifts left
# into a contiguous B:M range -- this is the "control de errores" cleanup
# described by the commit message.
#
# Delete right-to-left so the remaining column letters stay valid while we
# work.
$ws.Columns("I").Delete()
$ws.Columns("E:F").Delete()

# Reflect the new scroll position / active selection used while reviewing
# the corrected data.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("E5").Select()
